# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Estado de Cuenta" detail table (rows 16-21, columns C/D/E = N° Doc,
# Nombre, Periodo Mora) is re-sorted: it used to be grouped by worker
# (each worker's 1810 period followed by their 1809 period); now it is
# grouped by period (all three workers for period 1809 first, then all
# three workers for period 1810).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$doc1 = "73164755"
$name1 = "CAMILO YEPES CADENA"
$doc2 = "1143354313"
$name2 = "WILMER JAVIER DOMINGUEZ CASTRO"
$doc3 = "1002250621"
$name3 = "ANDREINA PAOLA ROMERO NARVAEZ"

$period1809 = "1809"
$period1810 = "1810"

# Row 16: Camilo Yepes Cadena - periodo 1809
$ws.Range("C16").Value = $doc1
$ws.Range("D16").Value = $name1
$ws.Range("E16").Value = $period1809

# Row 17: Wilmer Javier Dominguez Castro - periodo 1809
$ws.Range("C17").Value = $doc2
$ws.Range("D17").Value = $name2
$ws.Range("E17").Value = $period1809

# Row 18: Andreina Paola Romero Narvaez - periodo 1809
$ws.Range("C18").Value = $doc3
$ws.Range("D18").Value = $name3
$ws.Range("E18").Value = $period1809

# Row 19: Camilo Yepes Cadena - periodo 1810
$ws.Range("C19").Value = $doc1
$ws.Range("D19").Value = $name1
$ws.Range("E19").Value = $period1810

# Row 20: Wilmer Javier Dominguez Castro - periodo 1810
$ws.Range("C20").Value = $doc2
$ws.Range("D20").Value = $name2
$ws.Range("E20").Value = $period1810

# Row 21: Andreina Paola Romero Narvaez - periodo 1810
$ws.Range("C21").Value = $doc3
$ws.Range("D21").Value = $name3
$ws.Range("E21").Value = $period1810
